$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure B85 is stored as a real number (3) instead of text
$ws.Range("B85").Value = 3

# Add new row 86 with the annotation data
$ws.Range("A86").Value = "Ruilin"

# B86 keeps the politeness_score as text ("2"), not a number, so force
# text formatting before assigning, then restore the default cell style
# so no residual numeric formatting is left behind.
$ws.Range("B86").NumberFormat = "@"
$ws.Range("B86").Value = "2"
$ws.Range("B86").Style = "Normal"

$ws.Range("C86").Value = "No technical contribution."
$ws.Range("D86").Value = "CRT"
$ws.Range("E86").Value = "OTH"
$ws.Range("F86").Value = "e885cb01-c8a9-4c3f-b9a4-e5ab35292953"
$ws.Range("G86").Value = "SkwAEQbAb_annotated.xlsx"
$ws.Range("H86").Value = "No technical contribution."
